$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 169
$ws1.Range("F4").Value = 581
$ws1.Range("F6").Value = 476
$ws1.Range("F9").Value = 2446
$ws1.Range("G10").Value = 55
$ws1.Range("F12").Value = 166
$ws1.Range("F13").Value = 1484
$ws1.Range("F14").Value = 518
$ws1.Range("F17").Value = 225
$ws1.Range("F24").Value = 129
$ws1.Range("F26").Value = 1520
$ws1.Range("F27").Value = 15
$ws1.Range("F28").Value = 380
$ws1.Range("F29").Value = 363
$ws1.Range("F30").Value = 188
$ws1.Range("F31").Value = 291
$ws1.Range("F32").Value = 386

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 169
$ws4.Range("F4").Value = 581
$ws4.Range("F7").Value = 476
$ws4.Range("F10").Value = 2446
$ws4.Range("G11").Value = 55
$ws4.Range("F13").Value = 166
$ws4.Range("F14").Value = 1484
$ws4.Range("F15").Value = 518
$ws4.Range("F18").Value = 225
$ws4.Range("F25").Value = 129
$ws4.Range("F27").Value = 1520
$ws4.Range("F28").Value = 15
$ws4.Range("F29").Value = 380
$ws4.Range("F30").Value = 363
$ws4.Range("F31").Value = 188
$ws4.Range("F32").Value = 291
$ws4.Range("F33").Value = 386
